# Conformance fixture touch-up:
#  - sect1!E3 unit changes from "m" to "in" (new shared string)
#  - sect1 becomes the active / selected sheet (was sect2), with E4 selected
#  - sect2 keeps its own selection (A3) but is no longer the "tabSelected" sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("sect1")
$ws2 = $wb.Worksheets.Item("sect2")

# sect1: E3 changes units from "m" to "in"
$ws1.Range("E3").Value = "in"

# Make sect1 the active sheet and select E4 on it (mirrors the target
# sheetView: tabSelected="1" + <selection activeCell="E4" sqref="E4"/>).
[void]$ws1.Activate()
[void]$ws1.Range("E4").Select()

# sect2 retains its own cursor position (A3) but is no longer the active tab;
# activating sect1 above already clears sect2's tabSelected flag.
[void]$ws2.Range("A3").Select()

# Re-activate sect1 so it is left as the workbook's active sheet / tab.
[void]$ws1.Activate()
